# Platform Architecture slide rebuild (Slide 12 in the deck / "Slide 11" per
# the author's commit message). Replaces the old "UNIT ECONOMICS" slide
# content with a detailed 3-sided-marketplace architecture diagram.
#
# Colors are supplied as BGR-packed integers because PowerPoint's
# ColorFormat.RGB (and the underlying Windows COLORREF) stores 0x00BBGGRR,
# not 0x00RRGGBB:
#   3B82F6 (blue)   -> 0xF6823B
#   0A1628 (navy)   -> 0x28160A
#   475569 (slate)  -> 0x695547
#   FFFFFF (white)  -> 0xFFFFFF
#   A855F7 (purple) -> 0xF755A8
#   334155 (dk sl.) -> 0x554133
#   10B981 (green)  -> 0x81B910
#   C9A227 (gold)   -> 0x27A2C9
#
# Coordinates in the shape constructors are EMU/12700 (i.e. points), since
# Shapes.AddShape/AddTextbox take point values.

$p = $ppt.ActivePresentation
$sl = $p.Slides.Item(12)

function Pt($emu) {
    return $emu / 12700.0
}

# ---------------------------------------------------------------------
# 1) Header kicker bar -> "PLATFORM ARCHITECTURE"
# ---------------------------------------------------------------------
$sh2 = $sl.Shapes.Item(1)
$sh2.Left = Pt 1371600
$sh2.Top = Pt 411480
$sh2.Width = Pt 6400800
$sh2.Height = Pt 257175
$sh2.Fill.ForeColor.RGB = 0xF6823B
$sh2.Fill.Transparency = 0.85
$tr2 = $sh2.TextFrame.TextRange
$tr2.Text = "PLATFORM ARCHITECTURE"
$tr2.Font.Size = 10
$tr2.Font.Bold = 1
$tr2.Font.Color.RGB = 0xF6823B

# ---------------------------------------------------------------------
# 2) Headline -> "3-Sided Marketplace"
# ---------------------------------------------------------------------
$sh3 = $sl.Shapes.Item(2)
$sh3.Left = Pt 914400
$sh3.Top = Pt 720090
$sh3.Width = Pt 7315200
$sh3.Height = Pt 411480
$tr3 = $sh3.TextFrame.TextRange
$tr3.Text = "3-Sided Marketplace"
$tr3.Font.Size = 36
$tr3.Font.Bold = 1
$tr3.Font.Color.RGB = 0x28160A

# ---------------------------------------------------------------------
# 3) Subhead -> enterprise-grade platform description
# ---------------------------------------------------------------------
$sh4 = $sl.Shapes.Item(3)
$sh4.Left = Pt 1371600
$sh4.Top = Pt 1183005
$sh4.Width = Pt 6400800
$sh4.Height = Pt 257175
$tr4 = $sh4.TextFrame.TextRange
$tr4.Text = "Enterprise-grade platform connecting users, merchants, and payment partners"
$tr4.Font.Size = 12
$tr4.Font.Color.RGB = 0x695547

# ---------------------------------------------------------------------
# 4) Column 1 - User App card
# ---------------------------------------------------------------------
$shape3 = $sl.Shapes.AddShape(1, (Pt 731520), (Pt 1645920), (Pt 2560320), (Pt 1440180))
$shape3.Name = "Shape 3"
$shape3.Fill.ForeColor.RGB = 0xFFFFFF
$shape3.Line.ForeColor.RGB = 0xF6823B
$shape3.Line.Weight = 2
$shape3.Line.DashStyle = 1

$text4 = $sl.Shapes.AddShape(1, (Pt 731520), (Pt 1748790), (Pt 2560320), (Pt 257175))
$text4.Name = "Text 4"
$text4.Fill.Visible = 0
$text4.TextFrame.WordWrap = 1
$text4.TextFrame.VerticalAnchor = 3
$t4tr = $text4.TextFrame.TextRange
$t4tr.Text = "User App"
$t4tr.ParagraphFormat.Alignment = 2
$t4tr.ParagraphFormat.Bullet.Visible = 0
$t4tr.Font.Size = 14
$t4tr.Font.Bold = 1
$t4tr.Font.Color.RGB = 0x28160A

$text5 = $sl.Shapes.AddShape(1, (Pt 731520), (Pt 2057400), (Pt 2468880), (Pt 925830))
$text5.Name = "Text 5"
$text5.Fill.Visible = 0
$text5.TextFrame.WordWrap = 1
$text5.TextFrame.VerticalAnchor = 1
$t5tr = $text5.TextFrame.TextRange
$t5tr.Text = "• Search Engine: AI, voice, filters" + [char]13 + "• Wallet System: Dual coins" + [char]13 + "• Social: Referrals, leaderboards" + [char]13 + "• Gamification: Loyalty tiers"
$t5tr.ParagraphFormat.Bullet.Visible = 0
$t5tr.Font.Size = 8
$t5tr.Font.Color.RGB = 0x554133

# ---------------------------------------------------------------------
# 5) Column 2 - Merchant Hub card
# ---------------------------------------------------------------------
$shape6 = $sl.Shapes.AddShape(1, (Pt 3474720), (Pt 1645920), (Pt 2560320), (Pt 1440180))
$shape6.Name = "Shape 6"
$shape6.Fill.ForeColor.RGB = 0xFFFFFF
$shape6.Line.ForeColor.RGB = 0xF755A8
$shape6.Line.Weight = 2
$shape6.Line.DashStyle = 1

$text7 = $sl.Shapes.AddShape(1, (Pt 3474720), (Pt 1748790), (Pt 2560320), (Pt 257175))
$text7.Name = "Text 7"
$text7.Fill.Visible = 0
$text7.TextFrame.WordWrap = 1
$text7.TextFrame.VerticalAnchor = 3
$t7tr = $text7.TextFrame.TextRange
$t7tr.Text = "Merchant Hub"
$t7tr.ParagraphFormat.Alignment = 2
$t7tr.ParagraphFormat.Bullet.Visible = 0
$t7tr.Font.Size = 14
$t7tr.Font.Bold = 1
$t7tr.Font.Color.RGB = 0x28160A

$text8 = $sl.Shapes.AddShape(1, (Pt 3474720), (Pt 2057400), (Pt 2468880), (Pt 925830))
$text8.Name = "Text 8"
$text8.Fill.Visible = 0
$text8.TextFrame.WordWrap = 1
$text8.TextFrame.VerticalAnchor = 1
$t8tr = $text8.TextFrame.TextRange
$t8tr.Text = "• Analytics: Demographics, hours" + [char]13 + "• Campaign Manager: A/B testing" + [char]13 + "• QR/POS Integration" + [char]13 + "• Multi-Location Management"
$t8tr.ParagraphFormat.Bullet.Visible = 0
$t8tr.Font.Size = 8
$t8tr.Font.Color.RGB = 0x554133

# ---------------------------------------------------------------------
# 6) Column 3 - Admin Panel card
# ---------------------------------------------------------------------
$shape9 = $sl.Shapes.AddShape(1, (Pt 6217920), (Pt 1645920), (Pt 2560320), (Pt 1440180))
$shape9.Name = "Shape 9"
$shape9.Fill.ForeColor.RGB = 0xFFFFFF
$shape9.Line.ForeColor.RGB = 0x81B910
$shape9.Line.Weight = 2
$shape9.Line.DashStyle = 1

$text10 = $sl.Shapes.AddShape(1, (Pt 6217920), (Pt 1748790), (Pt 2560320), (Pt 257175))
$text10.Name = "Text 10"
$text10.Fill.Visible = 0
$text10.TextFrame.WordWrap = 1
$text10.TextFrame.VerticalAnchor = 3
$t10tr = $text10.TextFrame.TextRange
$t10tr.Text = "Admin Panel"
$t10tr.ParagraphFormat.Alignment = 2
$t10tr.ParagraphFormat.Bullet.Visible = 0
$t10tr.Font.Size = 14
$t10tr.Font.Bold = 1
$t10tr.Font.Color.RGB = 0x28160A

$text11 = $sl.Shapes.AddShape(1, (Pt 6217920), (Pt 2057400), (Pt 2468880), (Pt 925830))
$text11.Name = "Text 11"
$text11.Fill.Visible = 0
$text11.TextFrame.WordWrap = 1
$text11.TextFrame.VerticalAnchor = 1
$t11tr = $text11.TextFrame.TextRange
$t11tr.Text = "• Fraud Detection: 8-layer defense" + [char]13 + "• KYC Automation" + [char]13 + "• Content Moderation" + [char]13 + "• Financial Reporting"
$t11tr.ParagraphFormat.Bullet.Visible = 0
$t11tr.Font.Size = 8
$t11tr.Font.Color.RGB = 0x554133

# ---------------------------------------------------------------------
# 7) Intelligence Layer banner
# ---------------------------------------------------------------------
$shape12 = $sl.Shapes.AddShape(1, (Pt 731520), (Pt 3240405), (Pt 7680960), (Pt 822960))
$shape12.Name = "Shape 12"
$shape12.Fill.ForeColor.RGB = 0x28160A
$shape12.Line.ForeColor.RGB = 0x27A2C9
$shape12.Line.Weight = 2
$shape12.Line.DashStyle = 1

$text13 = $sl.Shapes.AddShape(1, (Pt 914400), (Pt 3343275), (Pt 7315200), (Pt 205740))
$text13.Name = "Text 13"
$text13.Fill.Visible = 0
$text13.TextFrame.WordWrap = 1
$text13.TextFrame.VerticalAnchor = 3
$t13tr = $text13.TextFrame.TextRange
$t13tr.Text = "Intelligence Layer"
$t13tr.ParagraphFormat.Bullet.Visible = 0
$t13tr.Font.Size = 14
$t13tr.Font.Bold = 1
$t13tr.Font.Color.RGB = 0x27A2C9

$text14 = $sl.Shapes.AddShape(1, (Pt 914400), (Pt 3600450), (Pt 7315200), (Pt 308610))
$text14.Name = "Text 14"
$text14.Fill.Visible = 0
$text14.TextFrame.WordWrap = 1
$text14.TextFrame.VerticalAnchor = 3
$t14tr = $text14.TextFrame.TextRange
$t14tr.Text = "Price Intelligence • Behavioral Insights • Predictive Analytics • Merchant ROI Tracking"
$t14tr.ParagraphFormat.Bullet.Visible = 0
$t14tr.Font.Size = 9
$t14tr.Font.Color.RGB = 0xFFFFFF

# ---------------------------------------------------------------------
# 8) Bottom banner -> "Not just an app..."
# ---------------------------------------------------------------------
$text15 = $sl.Shapes.AddShape(1, (Pt 1371600), (Pt 4371975), (Pt 6400800), (Pt 257175))
$text15.Name = "Text 15"
$text15.Fill.ForeColor.RGB = 0xF6823B
$text15.Fill.Transparency = 0.9
$text15.Line.Visible = 0
$text15.TextFrame.WordWrap = 1
$text15.TextFrame.VerticalAnchor = 3
$t15tr = $text15.TextFrame.TextRange
$t15tr.Text = "Not just an app • Built as a full-stack platform from Day 1"
$t15tr.ParagraphFormat.Alignment = 2
$t15tr.ParagraphFormat.Bullet.Visible = 0
$t15tr.Font.Size = 11
$t15tr.Font.Bold = 1
$t15tr.Font.Color.RGB = 0xF6823B

Write-Host "Platform Architecture slide rebuilt"
